$wb = $excel.ActiveWorkbook

# --- Sheet 2012 ---
$ws = $wb.Worksheets.Item("2012")
$ws.Range("E27").Value = -3264.8749993161441
$ws.Range("F27").Value = -3264.8749993161441
$ws.Range("E31").Value = -3264.8749993161591
$ws.Range("F31").Value = -3264.8749993161591
$ws.Range("A36").Value = "D.61"
$ws.Range("B36").Value = "Contribuciones sociales netas"
$ws.Range("E36").Value = 113976.0975413
$ws.Range("H36").Value = 137582.86030768359
$ws.Range("A37").Value = "D.73"
$ws.Range("B37").Value = "Transferencias corrientes dentro del gobierno general"
$ws.Range("E37").Value = 79482.446760280014
$ws.Range("F37").Value = 79482.446760280014
$ws.Range("H37").ClearContents()
$ws.Range("A38").Value = "D.76"
$ws.Range("B38").Value = "Transferencia del BCU al Gobierno"
$ws.Range("E38").Value = 926.5806651800001
$ws.Range("F38").ClearContents()

# --- Sheet 2016 ---
$ws = $wb.Worksheets.Item("2016")
$ws.Range("F27").Value = -10570.43407502865
$ws.Range("F31").Value = -10570.434075028639
$ws.Range("H34").Value = 55330.055140544151
$ws.Range("H35").Value = 55330.055140544697
$ws.Range("A36").Value = "D.76"
$ws.Range("B36").Value = "Transferencia del BCU al Gobierno"
$ws.Range("E36").Value = 1410.7265588499999
$ws.Range("G36").ClearContents()
$ws.Range("H36").ClearContents()
$ws.Range("A38").Value = "D.61"
$ws.Range("B38").Value = "Contribuciones sociales netas"
$ws.Range("E38").Value = 181558.9914336357
$ws.Range("G38").Value = 310.92989817709997
$ws.Range("H38").Value = 223683.25268700649

# --- Sheet 2017 ---
$ws = $wb.Worksheets.Item("2017")
$ws.Range("E27").Value = -25017.38878234318
$ws.Range("F27").Value = -25017.38878234318
$ws.Range("E31").Value = -25017.388782343201
$ws.Range("F31").Value = -25017.388782343201
$ws.Range("E34").Value = -63565.770979082881
$ws.Range("F34").Value = -63565.770979082881
$ws.Range("E35").Value = -63565.818181042087
$ws.Range("F35").Value = -63565.818181042087

# --- Sheet 2018 ---
$ws = $wb.Worksheets.Item("2018")
$ws.Range("E34").Value = -49773.892247343349
$ws.Range("F34").Value = -49773.892247342999
$ws.Range("H34").Value = 33283.712597372112
$ws.Range("E35").Value = -49773.892247342999
$ws.Range("F35").Value = -49773.892247343349
$ws.Range("H35").Value = 33283.712597372418
$ws.Range("A36").Value = "D.61"
$ws.Range("B36").Value = "Contribuciones sociales netas"
$ws.Range("E36").Value = 215925.528549476
$ws.Range("G36").Value = 381.05244255449912
$ws.Range("H36").Value = 271956.46735150547
$ws.Range("A38").Value = "D.76"
$ws.Range("B38").Value = "Transferencia del BCU al Gobierno"
$ws.Range("E38").Value = 1640.019891895929
$ws.Range("G38").ClearContents()
$ws.Range("H38").ClearContents()

# --- Sheet 2019 ---
$ws = $wb.Worksheets.Item("2019")
$ws.Range("F27").Value = -47487.012322160852
$ws.Range("F31").Value = -47487.012322160823
$ws.Range("D34").Value = -14516.140951382629
$ws.Range("H34").Value = 39248.472782472993
$ws.Range("D35").Value = -14516.15048177604
$ws.Range("H35").Value = 39248.472782472847
$ws.Range("A36").Value = "D.73"
$ws.Range("B36").Value = "Transferencias corrientes dentro del gobierno general"
$ws.Range("E36").Value = 200663.39336498949
$ws.Range("F36").Value = 200663.39336498949
$ws.Range("G36").ClearContents()
$ws.Range("H36").ClearContents()
$ws.Range("A37").Value = "D.61"
$ws.Range("B37").Value = "Contribuciones sociales netas"
$ws.Range("E37").Value = 239376.743428952
$ws.Range("F37").ClearContents()
$ws.Range("G37").Value = 402.67289956090139
$ws.Range("H37").Value = 298520.24294644757

# --- Sheet 2021 ---
$ws = $wb.Worksheets.Item("2021")
$ws.Range("G34").Value = 50492.474199902237
$ws.Range("H34").Value = 50492.474199901873
$ws.Range("G35").Value = 50492.474199901873
$ws.Range("H35").Value = 50492.474199902237
$ws.Range("A36").Value = "D.73"
$ws.Range("B36").Value = "Transferencias corrientes dentro del gobierno general"
$ws.Range("E36").Value = 264574.52162751998
$ws.Range("F36").Value = 264574.52162751998
$ws.Range("G36").ClearContents()
$ws.Range("H36").ClearContents()
$ws.Range("A37").Value = "D.61"
$ws.Range("B37").Value = "Contribuciones sociales netas"
$ws.Range("E37").Value = 272538.25067318691
$ws.Range("F37").ClearContents()
$ws.Range("G37").Value = 456.40441621169998
$ws.Range("H37").Value = 335741.85541967727
